$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

$updates = @{
    2  = 4
    3  = 8
    4  = 7
    6  = 7
    7  = 6
    8  = 9
    9  = 9
    10 = 7
    12 = 6
    14 = 4
    16 = 7
    17 = 7
    18 = 12
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
